# Generate Report for Handoff
#
# The 94b1938f-f6ce-419d-a8a1-4baeff0c8bd6 file just finished its de-de
# handback, so its "Latest Handback DateTime" needs to move from the stale
# placeholder timestamp to the real handback time (2016-03-24 18:50:50),
# both on the per-language "de-de" sheet and on the roll-up "Overview" sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 7 is the 94b1938f-f6ce-419d-a8a1-4baeff0c8bd6.md file;
# column D is "Latest Handoff Date" (the de-de handback datetime roll-up).
$wsOverview.Range("D7").Value = "2016-03-24 18:50:50"

# de-de sheet: row 7 is the same file; column E is "Latest Handback DateTime".
$wsDeDe.Range("E7").Value = "2016-03-24 18:50:50"
